$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update boundary values in columns C (start) and D (end) for several rows.
$ws.Range("C67").Value = 33.67
$ws.Range("D67").Value = 34.174999999999997

$ws.Range("C68").Value = 34.174999999999997
$ws.Range("D68").Value = 34.700000000000003

$ws.Range("C69").Value = 34.700000000000003
$ws.Range("D69").Value = 35.200000000000003

$ws.Range("C70").Value = 35.200000000000003
$ws.Range("D70").Value = 35.4

$ws.Range("D72").Value = 34.225000000000001

$ws.Range("C73").Value = 37.125
$ws.Range("D73").Value = 37.524999999999999

$ws.Range("C74").Value = 37.524999999999999
$ws.Range("D74").Value = 37.924999999999997

$ws.Range("C75").Value = 37.924999999999997
$ws.Range("D75").Value = 38.26

$ws.Range("C83").Value = 41.1
$ws.Range("D83").Value = 42.2

$ws.Range("C84").Value = 42.2
$ws.Range("D84").Value = 42.4

$ws.Range("C85").Value = 42.4
$ws.Range("D85").Value = 43.244999999999997

$ws.Range("C86").Value = 43.244999999999997
$ws.Range("D86").Value = 43.75

$ws.Range("C87").Value = 43.75
$ws.Range("D87").Value = 44.2

$ws.Range("C88").Value = 44.2
$ws.Range("D88").Value = 44.6

$ws.Range("C89").Value = 44.6
$ws.Range("D89").Value = 45.2

$ws.Range("C90").Value = 45.2
$ws.Range("D90").Value = 45.8

$ws.Range("C91").Value = 45.8
$ws.Range("D91").Value = 46.225000000000001

$ws.Range("C92").Value = 46.225000000000001
$ws.Range("D92").Value = 46.37

$ws.Range("C93").Value = 46.37
$ws.Range("D93").Value = 46.88

$ws.Range("C94").Value = 46.88
$ws.Range("D94").Value = 47.84

$ws.Range("C148").Value = 79.105000000000004

$ws.Range("D183").Value = 105.33

# Update sheet view: scroll position, zoom, and selection.
$excel.ActiveWindow.Zoom = 178
$ws.Range("E145").Select()
$excel.ActiveWindow.ScrollRow = 197
